$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 67. This shifts the existing rows 67-152
# down to rows 69-154 (matching the target dimension A1:T154), carrying
# over formatting (e.g. the date style on column D) the way Excel does
# when inserting rows.
$ws.Rows.Item(67).Insert()
$ws.Rows.Item(67).Insert()

# Populate the two newly inserted rows (67 and 68) with their data.
$ws.Range("A67").Value2 = 11
$ws.Range("B67").Value2 = "Vega Monumental Concepción"
$ws.Range("C67").Value2 = "Bíobío"
$ws.Range("D67").Value2 = 44482
$ws.Range("E67").Value2 = 8
$ws.Range("F67").Value2 = "Fruta"
$ws.Range("G67").Value2 = 100102
$ws.Range("H67").Value2 = "Cítricos"
$ws.Range("I67").Value2 = 100102005
$ws.Range("J67").Value2 = "Naranja"
$ws.Range("K67").Value2 = "Lane Late"
$ws.Range("L67").Value2 = "Primera"
$ws.Range("M67").Value2 = 100
$ws.Range("N67").Value2 = 7000
$ws.Range("O67").Value2 = 7500
$ws.Range("P67").Value2 = 7250
$ws.Range("Q67").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R67").Value2 = "Región de O'Higgins"
$ws.Range("S67").Value2 = 483
$ws.Range("T67").Value2 = 15

$ws.Range("A68").Value2 = 11
$ws.Range("B68").Value2 = "Vega Monumental Concepción"
$ws.Range("C68").Value2 = "Bíobío"
$ws.Range("D68").Value2 = 44482
$ws.Range("E68").Value2 = 8
$ws.Range("F68").Value2 = "Fruta"
$ws.Range("G68").Value2 = 100102
$ws.Range("H68").Value2 = "Cítricos"
$ws.Range("I68").Value2 = 100102005
$ws.Range("J68").Value2 = "Naranja"
$ws.Range("K68").Value2 = "Lane Late"
$ws.Range("L68").Value2 = "Segunda"
$ws.Range("M68").Value2 = 50
$ws.Range("N68").Value2 = 6500
$ws.Range("O68").Value2 = 6500
$ws.Range("P68").Value2 = 6500
$ws.Range("Q68").Value2 = "$/caja 15 kilos empedrada"
$ws.Range("R68").Value2 = "Región de O'Higgins"
$ws.Range("S68").Value2 = 433
$ws.Range("T68").Value2 = 15
